$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "95. Unique Binary Search Trees II" (row 96) is now finished -> mark column E as "Y".
$ws.Range("E96").Value = "Y"

# Re-apply the autofilter on the "Finished" column (D1:E184, field 2 = column E)
# so it only keeps rows with value "N" (the blank/empty criterion is no longer
# included). This naturally hides row 96 since its value no longer matches.
$ws.Range("D1:E184").AutoFilter(2, @("N"), 7)

# Re-applying the filter also re-evaluates every row in the range, which would
# incorrectly hide the still-blank "Finished" rows (they previously stayed
# visible only thanks to the separate "include blanks" flag). Restore those
# rows back to visible, matching rows still marked "Medium" with blank Finished.
for ($r = 2; $r -le 184; $r++) {
    $difficulty = $ws.Range("D$r").Text
    $finished = $ws.Range("E$r").Text
    if ($difficulty -eq "Medium" -and $finished -eq "" -and $r -ne 96) {
        $ws.Rows.Item($r).Hidden = $false
    }
}

# Update the active cell / selection on the sheet.
$ws.Range("G104").Select()
